# Updated symbol list data (Price/Volume columns) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.61"
$ws.Range("E2").Value = "'-1.14%"
$ws.Range("D3").Value = "'35.48"
$ws.Range("E3").Value = "'4.81%"
$ws.Range("D4").Value = "'5.022"
$ws.Range("E4").Value = "'-2.32%"
$ws.Range("D5").Value = "'0.07700"
$ws.Range("E5").Value = "'-1.68%"
$ws.Range("E6").Value = "'-9.08%"
$ws.Range("D7").Value = "'8.024"
$ws.Range("E7").Value = "'-0.41%"
$ws.Range("D8").Value = "'4.001"
$ws.Range("E8").Value = "'2.43%"
$ws.Range("D9").Value = "'0.9262"
$ws.Range("E9").Value = "'-0.96%"
$ws.Range("D10").Value = "'0.09190"
$ws.Range("E10").Value = "'-6.47%"
$ws.Range("D11").Value = "'0.1821"
$ws.Range("E11").Value = "'2.18%"
$ws.Range("D12").Value = "'0.08482"
$ws.Range("E12").Value = "'-0.49%"
$ws.Range("D13").Value = "'0.03562"
$ws.Range("E13").Value = "'6.32%"
$ws.Range("D14").Value = "'0.09962"
$ws.Range("E14").Value = "'0.29%"
$ws.Range("D15").Value = "'0.001479"
$ws.Range("E15").Value = "'0.03%"
$ws.Range("D16").Value = "'0.005759"
$ws.Range("E16").Value = "'0.19%"
$ws.Range("D17").Value = "'3.475"
$ws.Range("E17").Value = "'0.24%"
$ws.Range("E18").Value = "'1.04%"
$ws.Range("E19").Value = "'2.83%"
$ws.Range("D20").Value = "'0.1324"
$ws.Range("E20").Value = "'-1.32%"
$ws.Range("D21").Value = "'4.591"
$ws.Range("E21").Value = "'7.62%"
$ws.Range("E22").Value = "'-2.01%"
$ws.Range("D23").Value = "'0.04667"
$ws.Range("E23").Value = "'0.50%"
$ws.Range("D24").Value = "'0.001238"
$ws.Range("E24").Value = "'1.37%"
$ws.Range("D25").Value = "'0.004475"
$ws.Range("E25").Value = "'1.29%"
$ws.Range("D26").Value = "'0.0001309"
$ws.Range("E26").Value = "'1.03%"
$ws.Range("E27").Value = "'-20.21%"
$ws.Range("D39").Value = "'0.01724"
$ws.Range("E39").Value = "'-1.43%"
$ws.Range("D40").Value = "'0.04674"
$ws.Range("E40").Value = "'-3.04%"
$ws.Range("D41").Value = "'0.007931"
$ws.Range("E41").Value = "'1.71%"
$ws.Range("D42").Value = "'0.1399"
$ws.Range("E42").Value = "'-1.04%"
$ws.Range("D43").Value = "'0.007706"
$ws.Range("E43").Value = "'-21.41%"
$ws.Range("D44").Value = "'0.002237"
$ws.Range("E44").Value = "'7.60%"
$ws.Range("E45").Value = "'-1.68%"
$ws.Range("D46").Value = "'0.00006233"
$ws.Range("E46").Value = "'2.03%"
$ws.Range("D47").Value = "'0.00000000755"
$ws.Range("E47").Value = "'1.01%"
$ws.Range("D48").Value = "'3.347"
$ws.Range("E48").Value = "'19.82%"
$ws.Range("D49").Value = "'0.002704"
$ws.Range("E49").Value = "'35.58%"
$ws.Range("D50").Value = "'0.00002115"
$ws.Range("E50").Value = "'1.01%"
$ws.Range("D51").Value = "'0.0002014"
$ws.Range("E51").Value = "'1.01%"
